$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in R14's requirement text (row 15, column B) ---
$ws.Range("B15").Value = "A aplicação possuirá graficos para administração de recursos da maquina"

# --- Row 10 (R9) grew taller to fit its long wrapped text ---
$ws.Rows(10).RowHeight = 30

# --- Apply the existing alternating row style (rows 2/3 pattern: s=4/5/6 then s=1/2/3)
#     down across the 12 new rows (16-27) before filling in values ---
$ws.Range("A2:C3").Copy()
$ws.Range("A16:C27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- New backlog items: R15 (helpdesk) and R16 (OSHI API) ---
$ws.Range("A16").Value = "R15"
$ws.Range("B16").Value = "O sistema deve ter um helpdesk"
$ws.Range("C16").Value = "Importante"

$ws.Range("A17").Value = "R16"
$ws.Range("B17").Value = "O sistema deverá utilizar a API OSHI p/ capitar os dados"
$ws.Range("C17").Value = "Essencial"

# --- Reserved/placeholder backlog rows R17-R26 (ID only, rest pending) ---
$ws.Range("A18").Value = "R17"
$ws.Range("A19").Value = "R18"
$ws.Range("A20").Value = "R19"
$ws.Range("A21").Value = "R20"
$ws.Range("A22").Value = "R21"
$ws.Range("A23").Value = "R22"
$ws.Range("A24").Value = "R23"
$ws.Range("A25").Value = "R24"
$ws.Range("A26").Value = "R25"
$ws.Range("A27").Value = "R26"

# --- View state: zoomed to 80%, selection resting on the next empty row (B18) ---
$excel.ActiveWindow.Zoom = 80
[void]$ws.Range("B18").Select()
